$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "HTTP address" column (I) and "HTTP threads number" column (J)
# mirroring the existing ServerIp / NetThreadsNum columns (F / H).

$ws.Range("J2").Value = "HTTP线程数量"
$ws.Range("I2").Value = "HTTP地址"

$ws.Range("I3").Value = "string"
$ws.Range("J3").Value = "int32"

$ws.Range("I4").Value = "HttpIp"
$ws.Range("I5").Value = "127.0.0.1 10097"
$ws.Range("J4").Value = "HttpThreadsNum"

$ws.Range("J5").Value = 1

# Match column widths of the new columns to column H
$ws.Columns("I:J").ColumnWidth = $ws.Columns("H").ColumnWidth

# Update the last selected cell as recorded in the workbook
$ws.Range("N20").Select()
